$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 68
$ws.Range("D68").Value = 44176
$ws.Range("K68").Value = 11000
$ws.Range("L68").Value = 12000
$ws.Range("M68").Value = 11500
$ws.Range("P68").Value = 767

# Row 69
$ws.Range("I69").Value = "Segunda"
$ws.Range("J69").Value = 160
$ws.Range("K69").Value = 9000
$ws.Range("L69").Value = 10000
$ws.Range("M69").Value = 9500
$ws.Range("P69").Value = 633

# Row 70
$ws.Range("D70").Value = 44529
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 130
$ws.Range("K70").Value = 13000
$ws.Range("L70").Value = 14000
$ws.Range("M70").Value = 13500
$ws.Range("P70").Value = 900

# Row 71
$ws.Range("D71").Value = 44848
$ws.Range("J71").Value = 140
$ws.Range("K71").Value = 15000
$ws.Range("L71").Value = 16000
$ws.Range("M71").Value = 15500
$ws.Range("P71").Value = 1033

# Row 72
$ws.Range("I72").Value = "Segunda"
$ws.Range("J72").Value = 150
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 10000
$ws.Range("M72").Value = 9500
$ws.Range("P72").Value = 633

# Row 73
$ws.Range("D73").Value = 44267
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 160
$ws.Range("K73").Value = 12000
$ws.Range("L73").Value = 13000
$ws.Range("M73").Value = 12500
$ws.Range("P73").Value = 833

# Row 74
$ws.Range("D74").Value = 44452
$ws.Range("H74").Value = "Cristal"
$ws.Range("J74").Value = 120
$ws.Range("K74").Value = 37000
$ws.Range("L74").Value = 38000
$ws.Range("M74").Value = 37500
$ws.Range("P74").Value = 2500

# Row 75
$ws.Range("H75").Value = "Inferno"
$ws.Range("K75").Value = 38000
$ws.Range("L75").Value = 40000
$ws.Range("M75").Value = 39000
$ws.Range("P75").Value = 2600

# Row 76
$ws.Range("D76").Value = 44435
$ws.Range("K76").Value = 28000
$ws.Range("L76").Value = 30000
$ws.Range("M76").Value = 29000
$ws.Range("P76").Value = 1933

# Row 77
$ws.Range("D77").Value = 44610
$ws.Range("K77").Value = 14000
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = 14500
$ws.Range("P77").Value = 967

# Row 78
$ws.Range("D78").Value = 44389
$ws.Range("J78").Value = 200
$ws.Range("K78").Value = 21000
$ws.Range("L78").Value = 22000
$ws.Range("M78").Value = 21500
$ws.Range("P78").Value = 1433

# Row 79
$ws.Range("D79").Value = 44592
$ws.Range("J79").Value = 140
$ws.Range("K79").Value = 14000
$ws.Range("L79").Value = 15000
$ws.Range("M79").Value = 14500
$ws.Range("P79").Value = 967

# Row 80
$ws.Range("I80").Value = "Segunda"
$ws.Range("J80").Value = 130
$ws.Range("K80").Value = 10000
$ws.Range("L80").Value = 11000
$ws.Range("M80").Value = 10500
$ws.Range("P80").Value = 700

# Row 81
$ws.Range("D81").Value = 44806
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 150
$ws.Range("K81").Value = 14000
$ws.Range("L81").Value = 15000
$ws.Range("M81").Value = 14500
$ws.Range("P81").Value = 967

# Row 82
$ws.Range("I82").Value = "Segunda"
$ws.Range("J82").Value = 160
$ws.Range("K82").Value = 11000
$ws.Range("L82").Value = 12000
$ws.Range("M82").Value = 11500
$ws.Range("P82").Value = 767

# Row 83
$ws.Range("D83").Value = 44771
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 120

# Row 84
$ws.Range("D84").Value = 44526
$ws.Range("I84").Value = "Segunda"
$ws.Range("J84").Value = 150
$ws.Range("K84").Value = 12000
$ws.Range("L84").Value = 13000
$ws.Range("M84").Value = 12500
$ws.Range("P84").Value = 833

# Row 85
$ws.Range("D85").Value = 44876
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 170
$ws.Range("K85").Value = 14000
$ws.Range("L85").Value = 15000
$ws.Range("M85").Value = 14500
$ws.Range("P85").Value = 967
